# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) values are recalculated/regenerated from the
# source data pull (strikeouts per game, "K" instead of the previous
# "Strike#" metric). Update the per-row values in column G accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 0
    6  = 2
    7  = 1
    8  = 0
    9  = 0
    10 = 2
    11 = 2
    12 = 1
    13 = 3
    14 = 1
    15 = 2
    16 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
